$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Custom Shape?" value (P4): it was wrongly using the False text
# (mixed-up axis), should read True like the other rows' text flags.
# Copy the existing "True" text cell (G3) so the written value keeps its
# shared-string / text type (rather than being auto-coerced to a boolean).
$ws.Range("G3").Copy()
$ws.Range("P4").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Reflect the corrected viewport/selection: the window had scrolled to the
# wrong axis (column) before; move the active selection to P7 and scroll
# the sheet so column G is the leftmost visible column.
$ws.Range("P7").Select()
$excel.ActiveWindow.ScrollColumn = 7
